$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'28.013.30"
$ws.Range('E2').Value = '  +6.83%  '
$ws.Range('D3').Value = "'1.740.90"
$ws.Range('E3').Value = '  +5.15%  '
$ws.Range('D4').Value = "'1.002"
$ws.Range('E4').Value = '  -0.14%  '
$ws.Range('D5').Value = "'228.83"
$ws.Range('E5').Value = '  +4.36%  '
$ws.Range('D6').Value = "'0.5460"
$ws.Range('E6').Value = '  +4.16%  '
$ws.Range('D7').Value = "'1.002"
$ws.Range('E7').Value = '  -0.21%  '
$ws.Range('D8').Value = "'0.2774"
$ws.Range('E8').Value = '  +4.17%  '
$ws.Range('D9').Value = "'0.06746"
$ws.Range('E9').Value = '  +6.15%  '
$ws.Range('D10').Value = "'21.73"
$ws.Range('E10').Value = '  +5.04%  '
$ws.Range('D11').Value = "'0.07791"
$ws.Range('D12').Value = "'4.709"
$ws.Range('E12').Value = '  +2.41%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = "'1.776.91"
$ws.Range('E13').Value = '  +8.24%  '
$ws.Range('B14').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C14').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D14').Value = "'1.979.72"
$ws.Range('E14').Value = '  +5.08%  '
$ws.Range('D15').Value = "'0.6002"
$ws.Range('E15').Value = '  +6.70%  '
$ws.Range('D16').Value = "'0.0₅8434"
$ws.Range('E16').Value = '  +2.39%  '
$ws.Range('D17').Value = "'69.36"
$ws.Range('E17').Value = '  +6.04%  '
$ws.Range('D18').Value = "'27.993.06"
$ws.Range('E18').Value = '  +6.73%  '
$ws.Range('D19').Value = "'226.25"
$ws.Range('E19').Value = '  +17.75%  '
$ws.Range('D20').Value = "'4.852"
$ws.Range('E20').Value = '  +3.32%  '
$ws.Range('E21').Value = '  -0.20%  '
$ws.Range('D22').Value = "'10.94"
$ws.Range('E22').Value = '  +5.15%  '
$ws.Range('D23').Value = "'6.232"
$ws.Range('E23').Value = '  +3.90%  '
$ws.Range('D24').Value = "'1.003"
$ws.Range('E24').Value = '  -0.19%  '
$ws.Range('D25').Value = "'146.22"
$ws.Range('E25').Value = '  +1.76%  '
$ws.Range('E26').Value = '  +4.17%  '
$ws.Range('D27').Value = "'7.474"
$ws.Range('E27').Value = '  +2.78%  '
$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').Value = "'17.11"
$ws.Range('E28').Value = '  +7.28%  '
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').Value = "'1.655"
$ws.Range('E29').Value = '  +9.46%  '
$ws.Range('D30').Value = "'0.05699"
$ws.Range('E30').Value = '  +1.29%  '
$ws.Range('D31').Value = "'1.317"
$ws.Range('E31').Value = '  +3.02%  '
$ws.Range('D32').Value = "'3.721"
$ws.Range('E32').Value = '  +6.14%  '
$ws.Range('D33').Value = "'3.538"
$ws.Range('E33').Value = '  +5.43%  '
$ws.Range('E34').Value = '  +5.34%  '
$ws.Range('D35').Value = "'0.9842"
$ws.Range('E35').Value = '  +3.20%  '
$ws.Range('D36').Value = "'2.861"
$ws.Range('E36').Value = '  +2.04%  '
$ws.Range('D37').Value = "'2.453"
$ws.Range('E37').Value = '  +1.74%  '
$ws.Range('D38').Value = "'0.5966"
$ws.Range('E38').Value = '  +3.62%  '
$ws.Range('D39').Value = "'0.01676"
$ws.Range('E39').Value = '  +4.79%  '
$ws.Range('D40').Value = "'6.017"
$ws.Range('E40').Value = '  +0.23%  '
$ws.Range('D41').Value = "'0.8486"
$ws.Range('E41').Value = '  +0.75%  '
$ws.Range('D42').Value = "'1.047.63"
$ws.Range('E42').Value = '  +4.04%  '
$ws.Range('E43').Value = '  -0.15%  '
$ws.Range('E44').Value = '  +0.26%  '
$ws.Range('D45').Value = "'1.886.24"
$ws.Range('E45').Value = '  +5.09%  '
$ws.Range('D46').Value = "'0.0₈115"
$ws.Range('E46').Value = '  +11.43%  '
$ws.Range('D47').Value = "'60.15"
$ws.Range('E47').Value = '  +3.00%  '
$ws.Range('D48').Value = "'8.340"
$ws.Range('E48').Value = '  +3.76%  '
$ws.Range('B49').Value = 'Mantle'
$ws.Range('C49').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D49').Value = "'0.4426"
$ws.Range('E49').Value = '  +1.79%  '
$ws.Range('B50').Value = 'Frax'
$ws.Range('C50').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D50').Value = "'1.007"
$ws.Range('E50').Value = '  +0.08%  '
$ws.Range('D51').Value = "'0.05324"
$ws.Range('E51').Value = '  -0.49%  '
